$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (E column) -> new "Comment" category text, per the renamed/reordered
# change-log categories (Addon -> Add-On rename + reshuffled category list).
$map = @{
  2 = 'Alias Changed for Network'
  3 = 'Network Added to Add-On Package in Jan 2020'
  4 = 'Network Removed from Add-On Package in Jan 2020'
  5 = 'Network Added to Add-On Package in Jan 2020'
  6 = 'Network Removed from Add-On Package in Jan 2020'
  7 = 'Network Added to Add-On Package in Jan 2020'
  8 = 'Network Removed from Add-On Package in Jan 2020'
  9 = 'Network Added to Add-On Package in Jan 2020'
  10 = 'Network Removed from Base Service in Jan 2020'
  11 = 'Network Added to Add-On Package in Jan 2020'
  12 = 'Network Removed from Add-On Package'
  13 = 'Network Added to Add-On Package in Jan 2020'
  14 = 'Network Removed from Add-On Package'
  15 = 'Network Added to Add-On Package in Jan 2020'
  16 = 'Network Removed from Add-On Package'
  17 = 'Network Removed from Base Service in Jan 2020'
  18 = 'Name of Add-On Package Changed'
  19 = 'Name of Add-On Package Changed'
  20 = 'Name of Add-On Package Changed'
  21 = 'Name of Add-On Package Changed'
  22 = 'Name of Add-On Package Changed'
  23 = 'Name of Add-On Package Changed'
  24 = 'Name of Add-On Package Changed'
  25 = 'Name of Add-On Package Changed'
  26 = 'Name of Add-On Package Changed'
  27 = 'Name of Add-On Package Changed'
  28 = 'Name of Add-On Package Changed'
  29 = 'Name of Add-On Package Changed'
  30 = 'Name of Add-On Package Changed'
  31 = 'Name of Add-On Package Changed'
  32 = 'Name of Add-On Package Changed'
  33 = 'Name of Add-On Package Changed'
  34 = 'Name of Add-On Package Changed'
  35 = 'Name of Add-On Package Changed'
  36 = 'Name of Add-On Package Changed'
  37 = 'Name of Add-On Package Changed'
  38 = 'Name of Add-On Package Changed'
  39 = 'Name of Add-On Package Changed'
  40 = 'Name of Add-On Package Changed'
  41 = 'Name of Add-On Package Changed'
  42 = 'Network Added to Add-On Package in Jan 2020'
  43 = 'Network Added to Add-On Package in Jan 2020'
  44 = 'Network Added to Add-On Package in Jan 2020'
  45 = 'Network Added to Add-On Package in Jan 2020'
  46 = 'Network Added to Add-On Package in Jan 2020'
  47 = 'Network Added to Add-On Package in Jan 2020'
  48 = 'Network Added to Add-On Package in Jan 2020'
  49 = 'Network Added to Add-On Package in Jan 2020'
  50 = 'Network Added to Add-On Package in Jan 2020'
  51 = 'Network Added to Add-On Package in Jan 2020'
  52 = 'Network Added to Add-On Package in Jan 2020'
  53 = 'Network Added to Add-On Package in Jan 2020'
  54 = 'Network Added to Add-On Package in Jan 2020'
  55 = 'Network Added to Add-On Package in Jan 2020'
  56 = 'Network Moved from Base Service to Add-On Package'
  57 = 'Network Removed from Add-On Package'
  58 = 'Network Removed from Add-On Package'
  59 = 'Network Removed from Add-On Package'
  60 = 'Network Removed from Add-On Package'
  61 = 'Network Removed from Base Service in Jan 2020'
  62 = 'Network Removed from Base Service in Jan 2020'
  63 = 'Network Removed from Base Service in Jan 2020'
  64 = 'Network Removed from Base Service in Jan 2020'
  65 = 'Network Removed from Base Service in Jan 2020'
  66 = 'Network Removed from Base Service in Jan 2020'
  67 = 'New Network Added to Database in Jan 2020'
  68 = 'New Network Added to Database in Jan 2020'
  69 = 'New Network Added to Database in Jan 2020'
  70 = 'New Network Added to Database in Jan 2020'
  71 = 'New Network Added to Database in Jan 2020'
  72 = 'New Network Added to Database in Jan 2020'
  73 = 'New Network Added to Database in Jan 2020'
  74 = 'New Network Added to Database in Jan 2020'
  75 = 'New Network Added to Database in Jan 2020'
  76 = 'New Network Added to Database in Jan 2020'
  77 = 'New Network Added to Database in Jan 2020'
  78 = 'New Network Added to Database in Jan 2020'
  79 = 'Network Removed from Database in Jan 2020'
  80 = 'Network Added to Base Service in Jan 2020'
  81 = 'Network Added to Add-On Package in Jan 2020'
  82 = 'Network Moved from One Add-On Package to Another Add-On Package'
  83 = 'Network Moved from One Add-On Package to Another Add-On Package'
  84 = 'Network Removed from Add-On Package'
  85 = 'Network Added to Add-On Package in Jan 2020'
  86 = 'Network Moved from One Add-On Package to Another Add-On Package'
  87 = 'Network Moved from One Add-On Package to Another Add-On Package'
  88 = 'Network Removed from Add-On Package'
}

foreach ($row in $map.Keys) {
    $ws.Range("E$row").Value = $map[$row]
}

# Restore the selected range as recorded in the saved view state.
[void]$ws.Range("E3:E85").Select()

Write-Output "done"
